# Fixed naive component forecaster bug - Presentation state 11.02.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: drop stray first-year forecast cells (C2, E2) - bug produced
# spurious values here before the series actually starts.
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3: drop stray C3 value; E3 recalculated with corrected precision.
$ws.Range("C3").ClearContents()
$ws.Range("E3").Value = -0.6367039903685923

# Remaining rows: values recomputed by the fixed naive forecaster.
$ws.Range("C4").Value = -3.956152295564885

$ws.Range("C5").Value = 1.234995474941436
$ws.Range("E5").Value = -0.209816187756795

$ws.Range("C6").Value = 0.8993608108207818

$ws.Range("C8").Value = 0.02019328874802717

$ws.Range("E10").Value = -0.001769149545449711

$ws.Range("C12").Value = 0.0720185131838802

$ws.Range("E13").Value = -0.950584780912811

$ws.Range("C14").Value = -0.8017595264762423
$ws.Range("E14").Value = 0.0476740348578808

$ws.Range("C16").Value = 0.9704846793491706
$ws.Range("E16").Value = -0.8754609427830351

$ws.Range("C18").Value = 0.3928252664241683

$ws.Range("C19").Value = 0.3224026462283369
$ws.Range("E19").Value = -0.9749878381046684
